$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates to the cryptos table (Coin/Link/Price/Volume(1h) columns).
# Two pairs of rows (16/17, 20/21, 31/32) swapped their coin identity entirely
# (ranking order changed), so B/C/D/E are all rewritten for those rows.
$updates = @(
    @{ Row=2; D='28.460.55'; E='  -0.43%  ' }
    @{ Row=3; D='1.561.44'; E='  -1.85%  ' }
    @{ Row=4; E='  +0.70%  ' }
    @{ Row=5; D='211.62'; E='  -1.18%  ' }
    @{ Row=6; E='  -0.79%  ' }
    @{ Row=7; E='  +0.74%  ' }
    @{ Row=8; D='46.20'; E='  +4.04%  ' }
    @{ Row=9; D='24.01'; E='  +0.04%  ' }
    @{ Row=10; E='  -2.00%  ' }
    @{ Row=11; E='  -1.84%  ' }
    @{ Row=12; D='0.0882'; E='  -0.44%  ' }
    @{ Row=13; D='1.783.14'; E='  -1.73%  ' }
    @{ Row=14; D='1.563.82'; E='  -1.38%  ' }
    @{ Row=15; E='  -2.45%  ' }
    @{ Row=16; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='3.67'; E='  -2.91%  ' }
    @{ Row=17; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='28.445.69'; E='  -0.14%  ' }
    @{ Row=18; D='61.91'; E='  -3.31%  ' }
    @{ Row=19; D='226.50'; E='  -3.46%  ' }
    @{ Row=20; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.0₃0692'; E='  -2.53%  ' }
    @{ Row=21; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.32'; E='  -2.60%  ' }
    @{ Row=22; E='  +0.48%  ' }
    @{ Row=23; E='  -6.67%  ' }
    @{ Row=24; D='9.09'; E='  -3.47%  ' }
    @{ Row=25; D='2.07'; E='  +6.53%  ' }
    @{ Row=26; D='150.18'; E='  -0.83%  ' }
    @{ Row=27; D='14.96'; E='  -2.58%  ' }
    @{ Row=28; E='  -3.13%  ' }
    @{ Row=29; E='  -2.92%  ' }
    @{ Row=30; E='  +0.65%  ' }
    @{ Row=31; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.11'; E='  -3.50%  ' }
    @{ Row=32; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0464'; E='  -2.16%  ' }
    @{ Row=33; E='  -1.25%  ' }
    @{ Row=34; D='3.14'; E='  -0.58%  ' }
    @{ Row=35; D='1.395.13'; E='  -1.66%  ' }
    @{ Row=36; D='1.05'; E='  +0.01%  ' }
    @{ Row=37; E='  -4.02%  ' }
    @{ Row=38; E='  +1.73%  ' }
    @{ Row=39; E='  -0.19%  ' }
    @{ Row=40; E='  -1.39%  ' }
    @{ Row=41; E='  -1.81%  ' }
    @{ Row=42; E='  +0.61%  ' }
    @{ Row=43; D='0.787'; E='  -3.41%  ' }
    @{ Row=44; D='5.55'; E='  -2.26%  ' }
    @{ Row=45; E='  +1.12%  ' }
    @{ Row=46; D='0.980'; E='  +0.88%  ' }
    @{ Row=47; D='62.71'; E='  -2.90%  ' }
    @{ Row=48; D='1.696.04'; E='  -1.58%  ' }
    @{ Row=49; D='85.90'; E='  -1.91%  ' }
    @{ Row=50; D='0.0₆0102'; E='  +1.25%  ' }
    @{ Row=51; E='  -1.45%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B" + $u.Row).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $u.Row).Value = $u.C }

    # Price/Volume columns hold formatted text (thousands-dot prices, padded
    # percentages). Force text (NumberFormat "@") before assignment so Excel
    # doesn't reinterpret digit-and-dot strings as numbers (which would drop
    # trailing zeros, e.g. "46.20" -> 46.2), then restore the default style so
    # no stray formatting is left behind on the cell.
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
